$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in remaining columns for the existing gemma2b cpu script entry
$ws.Range("B2").Value = "google/gemma-2b"
$ws.Range("C2").Value = "cpu"
$ws.Range("D2").Value = 27.39
$ws.Range("E2").Value = 4.9400000000000004
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2.51"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "20.08"
$ws.Range("G2").Style = "Normal"

# Row 3: new llama_gguf_gemma2b.py entry
$ws.Range("A3").Value = "llama_gguf_gemma2b.py"
$ws.Range("B3").Value = "google/gemma-2b"
$ws.Range("C3").Value = "cpu"
$ws.Range("D3").Value = 115.31
$ws.Range("E3").Value = 95.26
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "0.00"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "20.08"
$ws.Range("G3").Style = "Normal"

# Remove the old rows that are no longer part of the table
$ws.Range("A4:G10").Clear()

# Resize columns A and B to fit the new (shorter) content
$ws.Columns.Item(1).ColumnWidth = 45.15
$ws.Columns.Item(2).ColumnWidth = 14.8

# Select the new used range, matching the saved selection state
$ws.Range("A1:G3").Select()
